$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.805.82'
$ws.Range('D2').Style = $origStyle
$ws.Range('E2').Value = '  +0.56%  '
$origStyle = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.137.60'
$ws.Range('D3').Style = $origStyle
$ws.Range('E3').Value = '  +1.38%  '
$ws.Range('E4').Value = '  -0.01%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '528.90'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +1.17%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.65'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  -1.30%  '
$ws.Range('E7').Value = '  -0.01%  '
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.137.12'
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  +1.39%  '
$ws.Range('E9').Value = '  +3.10%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.23'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('E11').Value = '  -0.14%  '
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.400'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  +3.90%  '
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.680.02'
$ws.Range('D13').Style = $origStyle
$ws.Range('E13').Value = '  +1.28%  '
$ws.Range('E14').Value = '  +2.71%  '
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.53'
$ws.Range('D15').Style = $origStyle
$ws.Range('E15').Value = '  -2.20%  '
$ws.Range('E16').Value = '  +0.64%  '
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '57.934.73'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  +0.60%  '
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.158.74'
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  +1.89%  '
$ws.Range('E19').Value = '  -1.86%  '
$ws.Range('E20').Value = '  -0.33%  '
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.97'
$ws.Range('D21').Style = $origStyle
$ws.Range('E21').Value = '  -1.10%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '351.30'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  +4.39%  '
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.78'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  -0.84%  '
$ws.Range('E24').Value = '  -0.05%  '
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '68.56'
$ws.Range('D25').Style = $origStyle
$ws.Range('E25').Value = '  +2.98%  '
$ws.Range('E26').Value = '  -0.42%  '
$ws.Range('E27').Value = '  +0.15%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.998'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  -0.44%  '
$ws.Range('E29').Value = '  +0.31%  '
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.54'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  +4.66%  '
$ws.Range('E31').Value = '  +0.01%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.14'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  -5.45%  '
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.87'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  +0.74%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '21.16'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  +1.28%  '
$ws.Range('E35').Value = '  -1.00%  '
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.01'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  +8.07%  '
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '158.09'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  +0.76%  '
$ws.Range('E38').Value = '  +1.53%  '
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '26.42'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  -2.50%  '
$ws.Range('E40').Value = '  -2.16%  '
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0672'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  +1.61%  '
$ws.Range('E42').Value = '  +7.55%  '
$ws.Range('E43').Value = '  +7.03%  '
$ws.Range('E44').Value = '  +3.10%  '
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.178.31'
$ws.Range('D45').Style = $origStyle
$ws.Range('E45').Value = '  +1.19%  '
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0272'
$ws.Range('D46').Style = $origStyle
$ws.Range('E46').Value = '  +4.47%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '36.66'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  -0.35%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.349.19'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  +2.07%  '
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  +0.02%  '
$ws.Range('E50').Value = '  -0.75%  '
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.04'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  +0.61%  '
